$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Row 4 "Name": value cleared
$meta.Range("B4").ClearContents()

# Row 5 "Title": value changed
$meta.Range("B5").Value = 'Mapping Métier/CDA/FHIR : "Acte"'

# Row 8 "Date": value changed
$meta.Range("B8").Value = "2026-01-07T15:20:53+00:00"

# --- Mapping Table 0 sheet ---
$t0 = $wb.Worksheets.Item("Mapping Table 0")

$t0.Range("D16").Value = "FRCDAActe.entryRelationship:frReferenceInterne"
$t0.Range("D17").Value = "FRCDAActe.entryRelationship:frReferenceInterne"
$t0.Range("D18").Value = "FRCDAActe.entryRelationship:frReferenceInterne"
$t0.Range("D19").Value = "FRCDAActe.entryRelationship:frSimpleObservationDifficulte"
$t0.Range("D20").Value = "FRCDAActe.entryRelationship:frSimpleObservationScores"

# --- Mapping Table 1 sheet ---
$t1 = $wb.Worksheets.Item("Mapping Table 1")

$t1.Range("A9").Value = "FRCDAActe.entryRelationship:frReferenceInterneDM"
$t1.Range("A10").Value = "FRCDAActe.entryRelationship:frSimpleObservationObservationsLiees"
$t1.Range("D11").Value = "FRProcedureActDocument.performer.actor.extension:Intervenant"
$t1.Range("D12").Value = "FRProcedureActDocument.performer.actor.extension:Informateur"
$t1.Range("D13").Value = "FRProcedureActDocument.performer.actor.extension:Participant"
$t1.Range("A14").Value = "FRCDAActe.entryRelationship:frReferenceInterneMotifActe"
$t1.Range("A15").Value = "FRCDAActe.entryRelationship:frReferenceInterneRencontreAssociee"
$t1.Range("D16").Value = "FRProcedureActDocument.recorder.extension:author"
$t1.Range("D17").Value = "FRProcedureActDocument.extension:priority"
